$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 23:28"

# Update country statistics rows (columns B,C,D,E,F,G,H)
$updates = @{
    4   = @{ B = 5522315; C = 46049; D = 2896717; E = 2453111; G = 952;  H = 172487 }
    5   = @{ B = 3317096; C = 38201;               E = 825562; G = 661;  H = 107232 }
    8   = @{ B = 583653;  C = 4513;  D = 466941;  E = 105035; G = 121;  H = 11677 }
    22  = @{ B = 224478;  C = 704;                 E = 12638 }
    23  = @{ B = 215521;  C = 3310;                E = 101267 }
    27  = @{ B = 121889;  C = 237;   D = 108218;  E = 4647;   G = 4;    H = 9024 }
    33  = @{ B = 92233;   C = 1153;                E = 23609 }
    53  = @{ B = 46430;   C = 378;   D = 42806;   E = 3454;   G = 2;    H = 170 }
    76  = @{ B = 16993;   C = 58;    D = 13759;   E = 3126 }
    87  = @{ B = 9638;    C = 33;    D = 8597;    E = 802;    G = 1;    H = 239 }
    102 = @{ B = 6693;    C = 17;    D = 5928;    E = 608 }
    106 = @{ B = 5176;    C = 104;   D = 2047;    E = 2999;   G = 2;    H = 130 }
    120 = @{ B = 3163;    C = 27;    D = 2302;    E = 827;    G = 1;    H = 34 }
    127 = @{ B = 2488;    C = 6;                  E = 1266 }
    138 = @{                        D = 1013;    E = 317 }
    141 = @{ B = 1593;    C = 78;    D = 408;     E = 1125;   G = 2;    H = 60 }
}

foreach ($rowNum in $updates.Keys) {
    $cols = $updates[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}
